# turkey_super-lig_2023-2024.xlsx edits:
#  1) Rows 46 and 47 had their match data (columns F:V) swapped - the
#     "Sivasspor vs Ankaragucu" match and the "Karagumruk vs Hatayspor"
#     match traded places (columns A:E - Indice/pais/torneio/temporada/
#     data_partida - stayed put).
#  2) A new match row (Adana Demirspor vs Besiktas) was appended as row 59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) swap F:V between row 46 and row 47 ---------------------------------
$row46 = $ws.Range("F46:V46").Value2
$row47 = $ws.Range("F47:V47").Value2
$ws.Range("F46:V46").Value2 = $row47
$ws.Range("F47:V47").Value2 = $row46

# --- 2) append new row 59 ---------------------------------------------------
# Copy cell formatting from the last existing data row (58) so the new row's
# styles (bold/bordered index cell, date-formatted match-date cell) match.
$ws.Range("A58").Copy()
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("E58").Copy()
$ws.Range("E59").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(59, 1).Value2  = 58
$ws.Cells.Item(59, 2).Value2  = "turkey"
$ws.Cells.Item(59, 3).Value2  = "super-lig"
$ws.Cells.Item(59, 4).Value2  = "2023-2024"
$ws.Cells.Item(59, 5).Value2  = 45196.79166666666
$ws.Cells.Item(59, 6).Value2  = "Adana Demirspor"
$ws.Cells.Item(59, 7).Value2  = 4
$ws.Cells.Item(59, 8).Value2  = "Besiktas"
$ws.Cells.Item(59, 9).Value2  = 2
$ws.Cells.Item(59, 10).Value2 = 2.73
$ws.Cells.Item(59, 11).Value2 = "21/08/2023 19:12"
$ws.Cells.Item(59, 12).Value2 = 3.25
$ws.Cells.Item(59, 13).Value2 = "27/09/2023 18:59"
$ws.Cells.Item(59, 14).Value2 = 3.88
$ws.Cells.Item(59, 15).Value2 = "21/08/2023 19:12"
$ws.Cells.Item(59, 16).Value2 = 4.04
$ws.Cells.Item(59, 17).Value2 = "27/09/2023 18:59"
$ws.Cells.Item(59, 18).Value2 = 2.45
$ws.Cells.Item(59, 19).Value2 = "21/08/2023 19:12"
$ws.Cells.Item(59, 20).Value2 = 2.11
$ws.Cells.Item(59, 21).Value2 = "27/09/2023 18:43"
$ws.Cells.Item(59, 22).Value2 = "https://www.betexplorer.com/football/turkey/super-lig/adanademirspor-besiktas/GfJ5g9O2/"

Write-Host "Applied row 46/47 swap and appended row 59"
